$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vscs")

# Insert a new row before row 18 (System info / Router ID), shifting rows 18+ down by one.
$ws.Rows.Item(18).Insert()

# Fix up formatting of the newly inserted row to match the row it displaced
# (row 19 now holds what used to be row 18 - "VM name" - with the correct per-column styles).
$ws.Range("A19:C19").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)

# Set the new cell's label.
$ws.Range("A18").Value = "Router ID"
